$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B gets narrower (44 -> 33 characters) ---
$ws.Columns.Item(2).ColumnWidth = 32.14

# --- Replace January dummy-data rows (2-29) with April data, and drop the
#     trailing January rows 30-35 that no longer exist in April (6 rows). ---
$ws.Rows("30:35").Delete()

$ws.Cells.Item(2, 1).Value = '''2025-04-01'
$ws.Cells.Item(2, 2).Value = 'Apple Pay Top-Up by *2180'
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 43.29
$ws.Cells.Item(2, 5).Value = 56.54

$ws.Cells.Item(3, 1).Value = '''2025-04-01'
$ws.Cells.Item(3, 2).Value = 'Apple Pay Top-Up by *1 573'
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 70.54

$ws.Cells.Item(4, 1).Value = '''2025-04-03'
$ws.Cells.Item(4, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 61.54

$ws.Cells.Item(5, 1).Value = '''2025-04-03'
$ws.Cells.Item(5, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(5, 3).Value = 9
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 52.54

$ws.Cells.Item(6, 1).Value = '''2025-04-04'
$ws.Cells.Item(6, 2).Value = 'Sainsbury''s'
$ws.Cells.Item(6, 3).Value = 3.75
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 48.79

$ws.Cells.Item(7, 1).Value = '''2025-04-05'
$ws.Cells.Item(7, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(7, 3).Value = 9
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 39.79

$ws.Cells.Item(8, 1).Value = '''2025-04-05'
$ws.Cells.Item(8, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(8, 3).Value = 9
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 30.79

$ws.Cells.Item(9, 1).Value = '''2025-04-06'
$ws.Cells.Item(9, 2).Value = 'To Ching Yin Siu'
$ws.Cells.Item(9, 3).Value = 10.15
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 20.64

$ws.Cells.Item(10, 1).Value = '''2025-04-11'
$ws.Cells.Item(10, 2).Value = 'Apple Pay Top-Up by *1 573'
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 30
$ws.Cells.Item(10, 5).Value = 50.64

$ws.Cells.Item(11, 1).Value = '''2025-04-12'
$ws.Cells.Item(11, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(11, 3).Value = 9
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 41.64

$ws.Cells.Item(12, 1).Value = '''2025-04-12'
$ws.Cells.Item(12, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(12, 3).Value = 9
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 32.64

$ws.Cells.Item(13, 1).Value = '''2025-04-12'
$ws.Cells.Item(13, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(13, 3).Value = 9
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 23.64

$ws.Cells.Item(14, 1).Value = '''2025-04-12'
$ws.Cells.Item(14, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 14.64

$ws.Cells.Item(15, 1).Value = '''2025-04-14'
$ws.Cells.Item(15, 2).Value = 'Apple Pay Top-Up by *2180'
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 50
$ws.Cells.Item(15, 5).Value = 64.64

$ws.Cells.Item(16, 1).Value = '''2025-04-15'
$ws.Cells.Item(16, 2).Value = 'Onerway'
$ws.Cells.Item(16, 3).Value = 31.26
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 33.38

$ws.Cells.Item(17, 1).Value = '''2025-04-17'
$ws.Cells.Item(17, 2).Value = 'Sainsbury''s'
$ws.Cells.Item(17, 3).Value = 3.69
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 29.69

$ws.Cells.Item(18, 1).Value = '''2025-04-17'
$ws.Cells.Item(18, 2).Value = 'KFC'
$ws.Cells.Item(18, 3).Value = 5.97
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 23.72

$ws.Cells.Item(19, 1).Value = '''2025-04-17'
$ws.Cells.Item(19, 2).Value = 'Www.cardiff.ac.uk'
$ws.Cells.Item(19, 3).Value = 9
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 14.72

$ws.Cells.Item(20, 1).Value = '''2025-04-21'
$ws.Cells.Item(20, 2).Value = 'Apple Pay Top-Up by *2180'
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 10
$ws.Cells.Item(20, 5).Value = 24.72

$ws.Cells.Item(21, 1).Value = '''2025-04-21'
$ws.Cells.Item(21, 2).Value = 'To Chung Ho Ling'
$ws.Cells.Item(21, 3).Value = 5.29
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 19.43

$ws.Cells.Item(22, 1).Value = '''2025-04-22'
$ws.Cells.Item(22, 2).Value = 'Apple Pay Top-Up by *2180'
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 9.19
$ws.Cells.Item(22, 5).Value = 28.62

$ws.Cells.Item(23, 1).Value = '''2025-04-23'
$ws.Cells.Item(23, 2).Value = 'Onerway'
$ws.Cells.Item(23, 3).Value = 28.62
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 10

$ws.Cells.Item(24, 1).Value = '''2025-04-24'
$ws.Cells.Item(24, 2).Value = 'From Ona*hungrypanda Co, London'
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 4.17
$ws.Cells.Item(24, 5).Value = 4.17

$ws.Cells.Item(25, 1).Value = '''2025-04-24'
$ws.Cells.Item(25, 2).Value = 'Apple Pay Top-Up by *1 573'
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 70
$ws.Cells.Item(25, 5).Value = 74.17

$ws.Cells.Item(26, 1).Value = '''2025-04-24'
$ws.Cells.Item(26, 2).Value = 'Reference: To Shuk Hei M'
$ws.Cells.Item(26, 3).Value = 14.9
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 59.27

$ws.Cells.Item(27, 1).Value = '''2025-04-24'
$ws.Cells.Item(27, 2).Value = 'To TS Kwong'
$ws.Cells.Item(27, 3).Value = 51.86
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 7.41

$ws.Cells.Item(28, 1).Value = '''2025-04-24'
$ws.Cells.Item(28, 2).Value = 'Payment from Kw'
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 4.5
$ws.Cells.Item(28, 5).Value = 1.91

$ws.Cells.Item(29, 1).Value = '''2025-04-27'
$ws.Cells.Item(29, 2).Value = 'giffgaff'
$ws.Cells.Item(29, 3).Value = 10
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 1.91

Write-Output "edit complete"
